# Add a new "reader" column to the site metadata sheet, inserted before the
# existing "array_type" column (which shifts from E -> F, F -> G, G -> H, H -> I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E; this shifts array_type/lat/lon/active right
# by one column and carries their formatting/widths along with them.
$ws.Columns("E:E").Insert()

# Header for the new column.
$ws.Range("E1").Value = "reader"
$ws.Range("E1").HorizontalAlignment = -4108  # xlCenter (matches header style of neighboring columns)

# Rows that are arrays/non litz_cord sites (no specific reader number) get "NA".
$naRows = @(2, 3, 22, 23, 24, 25, 26, 27, 28, 29, 30)
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 5).HorizontalAlignment = -4108  # xlCenter
}

# Rows 4-21 (the litz_cord sites) get sequential reader numbers 1-18.
for ($r = 4; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = $r - 3
    $ws.Cells.Item($r, 5).HorizontalAlignment = -4108  # xlCenter
}

# Match the column width Excel computed for the new "reader" column via AutoFit.
$ws.Columns("E:E").ColumnWidth = 6

# Update the active selection to reflect where the author ended up (F1).
$ws.Range("F1").Select()
